# Update attendee counts ("想去人数", column F) on the "展览" and "全部类型"
# sheets, which contain duplicated rows of event data.
$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 372
    $ws.Range("F3").Value = 106
    $ws.Range("F10").Value = 443
}
